# Add a new "Vue.js" benchmark section below the existing Jquery section,
# following the same layout/pattern used by the other library blocks
# (header row, Performance row, Memory row, Memory Leaks row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 44: section header ("Vue.js") ---
$ws.Range("B44").Value = "Vue.js"

# --- Row 45: Performance ---
$ws.Range("A45").Value = "Performance"
$ws.Range("B45").Value = 390
$ws.Range("C45").Value = 437
$ws.Range("D45").Value = 1070
$ws.Range("E45").Value = 1984
$ws.Range("F45").Value = 5839
$ws.Range("G45").Formula = "=SUM(B45:F45)"
$ws.Range("H45").Formula = "=AVERAGE(B45:F45)"
$ws.Range("J45").Value = 2

# --- Row 46: Memory ---
$ws.Range("A46").Value = "Memory"
$ws.Range("B46").Value = 11324
$ws.Range("C46").Value = 15451
$ws.Range("D46").Value = 28008
$ws.Range("E46").Value = 48566
$ws.Range("F46").Value = 112717
$ws.Range("G46").Value = 110181
$ws.Range("H46").Formula = "=AVERAGE(B46:G46)"
$ws.Range("I46").Value = "44833(memory taken at 0)"
$ws.Range("J46").Value = ""
$ws.Rows(46).RowHeight = 31.5

# --- Row 47: Memory Leaks ---
$ws.Range("A47").Value = "Memory Leaks"
$ws.Range("B47").Value = ""
$ws.Range("C47").Value = ""
$ws.Range("D47").Value = ""
$ws.Range("E47").Value = "200KB"
$ws.Range("F47").Value = ""
$ws.Range("G47").Value = ""

# Match the font/style used by the other data blocks (e.g. Backbone,
# Angular5 sections) which render with a plain black font instead of the
# theme-linked default font used by the section headers. Apply per-cell
# (rather than to the full rectangular block) so no stray blank cells are
# introduced outside the columns actually used by each row.
$ws.Range("A45:H45").Font.Color = 0
$ws.Range("J45").Font.Color = 0
$ws.Range("A46:J46").Font.Color = 0
$ws.Range("A47:G47").Font.Color = 0

# Leave the selection just past the newly added block, mirroring where the
# author's cursor ended up after typing in the new rows.
$ws.Range("A48").Select() | Out-Null
